# Rename header columns from "_old"/"_new" suffix to "_FV2310"/"_FV2404"
# and turn the data range A1:U65 into a real Excel Table ("Table1"),
# plus freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffixOld = "_FV2310"
$newSuffix = "_new"
$newSuffixNew = "_FV2404"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($text -like "*$oldSuffix") {
        $cell.Value = $text -replace [regex]::Escape($oldSuffix), $newSuffixOld
    } elseif ($text -like "*$newSuffix") {
        $cell.Value = $text -replace [regex]::Escape($newSuffix), $newSuffixNew
    }
}

# Convert the range into an Excel Table with header row, matching diff's table1.xml
$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1)
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
